$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab name Users -> Utilisateurs)
$ws.Name = "Utilisateurs"

# F1 header: refundBalance -> transactions
$ws.Range("F1").Value = "transactions"

# Row 2 edits
$ws.Range("D2").Value = '["jeudi","samedi","vendredi","dimanche"]'
$ws.Range("E2").Value = '[{"item":"Jardin","quantity":8},{"item":"Matelas","quantity":4}]'
$ws.Range("F2").Value = "[]"

# Row 3 edits
$ws.Range("C3").Value = $false
$ws.Range("F3").Value = '[{"id":"1745428649532","payer":"Ju","amount":666,"description":"d","repayments":[{"userId":"Mila","amount":222,"paid":false},{"userId":"Ju","amount":222,"paid":false},{"userId":"Louise","amount":222,"paid":false}]}]'

# New row 4
$ws.Range("A4").Value = "'1745427357665"
$ws.Range("B4").Value = "Louise"
$ws.Range("C4").Value = "'"
$ws.Range("D4").Value = "[]"
$ws.Range("E4").Value = "[]"
$ws.Range("F4").Value = "[]"

# Extend the "numbers stored as text" ignored-error range to cover the new row
$ws.Range("A1:F4").IgnoredErrors.NumberAsText = $true
